$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# Update the confidential footer date string in cell A40
$ws.Range("A40").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-06-10 for illustrative purposes only and are subject to change."

# Updated Weight (D) and Percent Change (E) values for rows 2-37
$ws.Range("D2").Value = 0.0312167796058197
$ws.Range("E2").Value = 0.01964912280701747

$ws.Range("D3").Value = 0.03020047392391594
$ws.Range("E3").Value = 0.04049117898500043

$ws.Range("D4").Value = 0.03127936964011708
$ws.Range("E4").Value = 0.02001000500250139

$ws.Range("D5").Value = 0.06417727844839305
$ws.Range("E5").Value = 0.02087682672233826

$ws.Range("D6").Value = 0.03034834287994351
$ws.Range("E6").Value = 0.002706883217324041

$ws.Range("D7").Value = 0.01559548060833601
$ws.Range("E7").Value = 0.01000827752276323

$ws.Range("D8").Value = 0.03279404847011375
$ws.Range("E8").Value = 0.02137608550434189

$ws.Range("D9").Value = 0.03132240028869653
$ws.Range("E9").Value = -0.00573248407643312

$ws.Range("D10").Value = 0.04751072587814307
$ws.Range("E10").Value = 0.003005290134002969

$ws.Range("D11").Value = 0.02920685712944499
$ws.Range("E11").Value = -0.001171948246765475

$ws.Range("D12").Value = 0.01558100666290474
$ws.Range("E12").Value = -0.01107205623901575

$ws.Range("D13").Value = 0.01762691840900044
$ws.Range("E13").Value = -0.01780958721704384

$ws.Range("D14").Value = 0.01408451806162576
$ws.Range("E14").Value = 0.01360940993486914

$ws.Range("D15").Value = 0.007219173674631574
$ws.Range("E15").Value = -0.006908884012029604

$ws.Range("D16").Value = 0.007834707543174899
$ws.Range("E16").Value = -0.02156980227681238

$ws.Range("D17").Value = 0.03232442761902619
$ws.Range("E17").Value = 0.01736625863018326

$ws.Range("D18").Value = 0.02897038415611519
$ws.Range("E18").Value = 0.01864092090605274

$ws.Range("D19").Value = 0.03225460061201318
$ws.Range("E19").Value = 0.01095169369216387

$ws.Range("D20").Value = 0.03229743566673544
$ws.Range("E20").Value = 0.006691900075700197

$ws.Range("D21").Value = 0.04873025357765614
$ws.Range("E21").Value = 0.01212169864333301

$ws.Range("D22").Value = 0.02994776666044026
$ws.Range("E22").Value = -0.0230942055488792

$ws.Range("D23").Value = 0.03005690803274632
$ws.Range("E23").Value = 0.003156113750244049

$ws.Range("D24").Value = 0.02914778778457684
$ws.Range("E24").Value = -0.008334339896122755

$ws.Range("D25").Value = 0.01451932320613539
$ws.Range("E25").Value = -0.02327837051406412

$ws.Range("D26").Value = 0.01486552433334279
$ws.Range("E26").Value = -0.003947264545669849

$ws.Range("D27").Value = 0.03241713910732919
$ws.Range("E27").Value = -0.0008628127696290733

$ws.Range("D28").Value = 0.02981632758841576
$ws.Range("E28").Value = 0.02602991340855398

$ws.Range("D29").Value = 0.02976038774526247
$ws.Range("E29").Value = 0.01439331203911842

$ws.Range("D30").Value = 0.02836580354357391
$ws.Range("E30").Value = -0.006523058252427161

$ws.Range("D31").Value = 0.02716133657106365
$ws.Range("E31").Value = 0.003845433727478254

$ws.Range("D32").Value = 0.02857626253389886
$ws.Range("E32").Value = 0.01026694045174548

$ws.Range("D33").Value = 0.02978288103883809
$ws.Range("E33").Value = 0.009476649876205778

$ws.Range("D34").Value = 0.03138107844585033
$ws.Range("E34").Value = 0.0009723261032161812

$ws.Range("D35").Value = 0.03187006308879863
$ws.Range("E35").Value = -0.01558856020621091

$ws.Range("D36").Value = 0.03175622746392026
$ws.Range("E36").Value = -0.001207208760886425

$ws.Range("D37").Value = 0.9999999999999999
$ws.Range("E37").Value = 0.006348585416326324

$ws.Protect()
